# Ajustes en vista de impresion de rotulos: 3 columnas, altura flexible
# Rebuild the "Servicio/Dieta" breakdown rows (A2:E32) for fecha 2025-06-xx (serial 45822),
# replacing the old single-column (Desayuno only, serial 45825) report with the new
# 3-column breakdown across Almuerzo + Desayuno, plus newly appearing diet categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# r=row, v=Fecha (date serial), servicio, dieta, cantidad, valorTotal
$data = @(
    @(45822, 'Almuerzo', 'Blanda', 100, 1410000),
    @(45822, 'Desayuno', '3.0 Onzas', 2, 0),
    @(45822, 'Desayuno', '3.5 Onzas', 3, 0),
    @(45822, 'Desayuno', 'Alta en Fibra', 3, 0),
    @(45822, 'Desayuno', 'Astringente', 29, 0),
    @(45822, 'Desayuno', 'Blanda', 205, 0),
    @(45822, 'Desayuno', 'Coronaria', 207, 0),
    @(45822, 'Desayuno', 'Hepatica', 3, 0),
    @(45822, 'Desayuno', 'Hipercalorica', 10, 0),
    @(45822, 'Desayuno', 'Hiperproteica', 22, 0),
    @(45822, 'Desayuno', 'Hipo Grasa', 29, 0),
    @(45822, 'Desayuno', 'Hipoglucida', 381, 0),
    @(45822, 'Desayuno', 'Hipograsa', 88, 0),
    @(45822, 'Desayuno', 'Hiposodica', 382, 0),
    @(45822, 'Desayuno', 'Liquida Clara', 61, 0),
    @(45822, 'Desayuno', 'Liquida Total', 70, 0),
    @(45822, 'Desayuno', 'Liquida Total 140 Cc', 4, 0),
    @(45822, 'Desayuno', 'Liquida Total Miel 140 Cc', 28, 0),
    @(45822, 'Desayuno', 'Liquida Total Nectar', 112, 0),
    @(45822, 'Desayuno', 'Liquida Total Nectar 140 Cc', 21, 0),
    @(45822, 'Desayuno', 'Liquida total Miel', 29, 0),
    @(45822, 'Desayuno', 'Nada Via Oral', 413, 0),
    @(45822, 'Desayuno', 'Normal', 646, 0),
    @(45822, 'Desayuno', 'Pequena Semiblanda', 52, 0),
    @(45822, 'Desayuno', 'Renal Dialisis', 106, 0),
    @(45822, 'Desayuno', 'Renal PRE Dialisis', 14, 0),
    @(45822, 'Desayuno', 'Renal SIN Dialisis', 47, 0),
    @(45822, 'Desayuno', 'Semiblanda', 299, 0),
    @(45822, 'Desayuno', 'Semiblanda Pequena', 11, 0),
    @(45822, 'Desayuno', 'Todo Pure', 28, 0),
    @(45822, 'Desayuno', 'nan', 100, 0)
)

$row = 2
foreach ($item in $data) {
    $fecha = $item[0]
    $servicio = $item[1]
    $dieta = $item[2]
    $cantidad = $item[3]
    $valorTotal = $item[4]

    $ws.Cells.Item($row, 1).Value = $fecha
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 2).Value = $servicio
    $ws.Cells.Item($row, 3).Value = $dieta
    $ws.Cells.Item($row, 4).Value = $cantidad
    $ws.Cells.Item($row, 5).Value = $valorTotal

    $row++
}
